# New weekly report: insert a new row for the latest week right after the
# header/most-recent rows (row 4), pushing the existing rows 4-18 down to
# 5-19 (the table is sorted most-recent-first with the header in row 1 and
# the two prior rows in 2-3), and populate the new row with this week's
# data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 4:18 down to 5:19, leaving a blank row 4 for the new entry.
$ws.Rows("4:4").Insert()

# Populate the new row 4 with this week's record.
$ws.Range("A4").Value = 7
$ws.Range("B4").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C4").Value = "Ñuble"
$ws.Range("D4").Value = 44959
$ws.Range("E4").Value = 16
$ws.Range("F4").Value = 100114007
$ws.Range("G4").Value = "Jengibre"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 30
$ws.Range("K4").Value = 19000
$ws.Range("L4").Value = 19000
$ws.Range("M4").Value = 19000
$ws.Range("N4").Value = "$/caja 13 kilos"
$ws.Range("O4").Value = "Perú"
$ws.Range("P4").Value = 1462
$ws.Range("Q4").Value = 13
$ws.Range("R4").Value = "Hortaliza"
